$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 121, pushing existing rows 121:137 down to 122:138
$ws.Rows("121:121").Insert()

# Populate the newly inserted row 121 with this week's data (same template as the
# surrounding Mango / Terminal Hortofruticola Agro Chillan rows, new date + prices)
$ws.Range("A121").Value = 7
$ws.Range("B121").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C121").Value = "Ñuble"
$ws.Range("D121").Value = 45077
$ws.Range("D121").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E121").Value = 16
$ws.Range("F121").Value = "Fruta"
$ws.Range("G121").Value = 100108
$ws.Range("H121").Value = "Tropicales y subtropicales"
$ws.Range("I121").Value = 100108002
$ws.Range("J121").Value = "Mango"
$ws.Range("K121").Value = "Sin especificar"
$ws.Range("L121").Value = "Primera"
$ws.Range("M121").Value = 50
$ws.Range("N121").Value = 8000
$ws.Range("O121").Value = 9000
$ws.Range("P121").Value = 8400
$ws.Range("Q121").Value = "$/bandeja 4 kilos"
$ws.Range("R121").Value = "Perú"
$ws.Range("S121").Value = 2100
$ws.Range("T121").Value = 4
